# Update average_county_temperature (column AD) with NOAA-sourced values,
# keyed by facility_id (column H). Facilities not present in the map keep
# their existing value (10) unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tempByFacility = @{
    1001726 = 13.75752314814816
    1004528 = 13.75752314814816
    1004529 = 3.38888888888889
    1004861 = 12.93898809523811
    1013270 = 19.79629629629628
    1013627 = 13.75752314814816
    1013656 = 3.38888888888889
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $facilityId = $ws.Cells.Item($r, 8).Value()
    if ($tempByFacility.ContainsKey($facilityId)) {
        $ws.Cells.Item($r, 30).Value = $tempByFacility[$facilityId]
    }
}
